# Fruta / hortaliza, semanal
# Insert a new weekly record at row 395 (pushes existing rows 395-421 down
# to 396-422) for "Ají" - "Inferno" variety at "Vega Central Mapocho de
# Santiago".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 395..421 down to 396..422, leaving a blank row 395 to fill in.
$ws.Rows.Item(395).Insert()

$ws.Range("A395").Value = 9
$ws.Range("B395").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C395").Value = 'Metropolitana'
$ws.Range("D395").Value = 45106
$ws.Range("E395").Value = 13
$ws.Range("F395").Value = 100112021
$ws.Range("G395").Value = 'Ají'
$ws.Range("H395").Value = 'Inferno'
$ws.Range("I395").Value = 'Primera'
$ws.Range("J395").Value = 70
$ws.Range("K395").Value = 11000
$ws.Range("L395").Value = 13000
$ws.Range("M395").Value = 12000
$ws.Range("N395").Value = '$/caja 10 kilos'
$ws.Range("O395").Value = 'Región de Arica y Parinacota'
$ws.Range("P395").Value = 1200
$ws.Range("Q395").Value = 10
$ws.Range("R395").Value = 'Hortaliza'
